$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-02 Monday", "2026-02-03 Tuesday"),
    @("36×53=", "21×75="),
    @("73×51=", "14×27="),
    @("73×93=", "98×34="),
    @("27×81=", "38×16="),
    @("41×21=", "79×86="),
    @("87×79=", "58×42="),
    @("26×80=", "64×91="),
    @("23×75=", "72×44="),
    @("97×94=", "75×14="),
    @("27×47=", "68×24="),
    @("97×11=", "57×19="),
    @("31×88=", "36×49="),
    @("66×86=", "67×23="),
    @("21×25=", "43×58="),
    @("46×36=", "95×65="),
    @("92×12=", "22×55="),
    @("43×12=", "69×82="),
    @("75×44=", "48×21="),
    @("38×15=", "41×33="),
    @("39×43=", "76×37="),
    @("91×56=", "36×50="),
    @("69×26=", "44×72="),
    @("43×66=", "34×53="),
    @("45×56=", "69×11="),
    @("93×93=", "85×60=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
